$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the erroneous "RBC / Deposit / 158.99" row (row 18) and shift rows up.
$ws.Rows.Item(18).Delete()

# Reflect the selection left behind after the row delete/cleanup operation.
$ws.Range("A18:XFD18").Select()
